# Tambahkan data awal pengguna dan kendaraan
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: new user "Sinta Maharani" (NIK only, no Plat yet)
# Leading "'" forces the long NIK digit-string to be stored as text so the
# 16-digit number keeps its exact digits instead of being rounded as a
# floating point number.
$ws.Range("A7").Value = "'1245367800112234"
$ws.Range("C7").Value = "Sinta Maharani"

# Row 8: new user "Dinda" (NIK only, no Plat yet)
$ws.Range("A8").Value = "'9801234567819235"
$ws.Range("C8").Value = "Dinda"
